$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: add PriceChange (X7) and UpDown (Y7) values ---
$ws.Range("X7").Value = 0.11999500000000296
$ws.Range("Y7").Value = "Up"

# --- Row 8: new data row ---
$ws.Range("A8").Value = 42647.887002314812
$ws.Range("B8").Value = -3
$ws.Range("C8").Value = "Neutral"
$ws.Range("D8").Value = 28
$ws.Range("E8").Value = 23176
$ws.Range("F8").Value = 2691
$ws.Range("G8").Value = 63
$ws.Range("H8").Value = 34
$ws.Range("I8").Value = 73
$ws.Range("J8").Value = 26
$ws.Range("K8").Value = 27358
$ws.Range("L8").Value = 332
$ws.Range("M8").Value = 184
$ws.Range("N8").Value = 73
$ws.Range("O8").Value = 26
$ws.Range("P8").Value = "Bag"
$ws.Range("Q8").Value = 42.459412013272512
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = -0.0112
$ws.Range("T8").Value = -0.0367
$ws.Range("U8").Value = 14.56
$ws.Range("V8").Value = "N/A"
$ws.Range("W8").Value = -2

# S8/T8 need the same percentage number format as the rest of the column
# (A8 already inherits the date style from the column's default style)
$ws.Range("S7:T7").Copy() | Out-Null
$ws.Range("S8:T8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column C is a touch wider now that "Neutral" is present (bestFit recalculated by Excel)
$ws.Columns.Item(3).ColumnWidth = 5.666666666666667
